# Update odds values in the "Jogos da Semana" worksheet to reflect the
# latest FlashScore data refresh.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6
$ws.Range("Q6").Value = 1.8
$ws.Range("R6").Value = 2

# Row 15
$ws.Range("Q15").Value = 2.05
$ws.Range("R15").Value = 1.8

# Row 16
$ws.Range("S16").Value = 1.41
$ws.Range("T16").Value = 2.62

# Row 17
$ws.Range("G17").Value = 1.53
$ws.Range("I17").Value = 5.75
$ws.Range("M17").Value = 1.03
$ws.Range("N17").Value = 15
$ws.Range("O17").Value = 1.2
$ws.Range("P17").Value = 4.33
$ws.Range("S17").Value = 1.3
$ws.Range("W17").Value = 7.5
$ws.Range("X17").Value = 7.5
$ws.Range("Z17").Value = 11
$ws.Range("AG17").Value = 251
$ws.Range("AU17").Value = 8.5
$ws.Range("AW17").Value = 7.5
$ws.Range("BA17").Value = 126
